# Update de lista de itens
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column layout tweaks -------------------------------------------------
# Column A ("Nome / Descricao") narrower: 34 -> 27 characters.
$ws.Columns.Item(1).ColumnWidth = 26.2857142857143

# Column F ("Atribuido a") is now hidden from view.
$ws.Columns.Item(6).Hidden = $true

# --- Row height tweaks (wrap height no longer needed / now needed) -------
$ws.Rows.Item(26).RowHeight = 15
$ws.Rows.Item(28).RowHeight = 15
$ws.Rows.Item(29).RowHeight = 15
$ws.Rows.Item(31).RowHeight = 15
$ws.Rows.Item(34).RowHeight = 30

# --- Status / hours updates on the work-item list -------------------------
$ws.Range("D34").Value = "Finalizado"
$ws.Range("D35").Value = "Finalizado"
$ws.Range("H36").Value = 42
$ws.Range("D37").Value = "Finalizado"
$ws.Range("H37").Value = 10
$ws.Range("D38").Value = "Iniciado"
$ws.Range("H38").Value = 8

# --- Leave the selection where the author left it before saving ----------
$ws.Range("D36").Select() | Out-Null
